# FWHM data run for sg_rr_100_025 2023-12-11 14-23-14
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 59

$ws.Cells.Item($row, 1).Value  = "sg_rr_100_025 2023-12-11 14-23-14.csv"
$ws.Cells.Item($row, 2).Value  = 0.01
$ws.Cells.Item($row, 3).Value  = 1000
$ws.Cells.Item($row, 4).Value  = 5001
$ws.Cells.Item($row, 5).Value  = 1530
$ws.Cells.Item($row, 6).Value  = 1570
$ws.Cells.Item($row, 7).Value  = 0.5
$ws.Cells.Item($row, 8).Value  = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = 0.98282051282051597
$ws.Cells.Item($row, 11).Value = 0.0055097596875867197
$ws.Cells.Item($row, 12).Value = "yes(although possible end peak not found)"
$ws.Cells.Item($row, 13).Value = 0.14513350341842499
$ws.Cells.Item($row, 14).Value = 0.0036926757851065001
$ws.Cells.Item($row, 15).Value = "reduced approx fsr a bit, to see if this had any affect on fsr calculation as above, half the approx fsr was quite close to actual calculated fsr."

# Scroll the sheet view up a bit and keep the active cell selection as before
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("A59").Select()
